# Revert "Predicting PM10 hourly levels / Organizing Data Gathering and
# Exploration code": drop the extra "pm10_limits" sheet (and the shared
# strings that only it used), and rename "arima_graph" back to "Sheet1" -
# restoring the chart's series references to point at the new sheet name.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the "pm10_limits" worksheet that was added by the reverted commit.
# sharedStrings.xml entries that were only referenced by this sheet are
# dropped automatically once the sheet is gone.
$wb.Worksheets.Item("pm10_limits").Delete()

# Rename the remaining sheet back to "Sheet1".
$ws = $wb.Worksheets.Item("arima_graph")
$ws.Name = "Sheet1"

# The line chart's two series still reference the old sheet name in their
# cached SERIES() formulas - repoint them at the renamed sheet.
$chart = $ws.ChartObjects(1).Chart
$series1 = $chart.SeriesCollection(1)
$series1.Formula = "=SERIES(Sheet1!`$B`$1,,Sheet1!`$B`$2:`$B`$13,1)"
$series2 = $chart.SeriesCollection(2)
$series2.Formula = "=SERIES(Sheet1!`$C`$1,,Sheet1!`$C`$2:`$C`$13,2)"
